# Update Name of Algo
# Applies corrected values to column D for specific rows in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -7.712000000000001
    14 = -7.542
    21 = -8.1
    23 = -7.874
    25 = -8.229000000000001
    26 = -8.238
    29 = -7.37
    53 = -7.507000000000001
    57 = -7.904000000000001
    59 = -8.061
    69 = -7.221000000000001
    79 = -7.885
    83 = -8.199
    91 = -6.787999999999999
    93 = -7.595000000000001
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
